# Lecture partielle de l'EDT M1 MIAGE.
# Update the day-of-week labels and the corresponding date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Day-of-week text labels (shared strings "jeudi" -> "lundi", "samedi" -> "mercredi")
$ws.Range("B2").Value = "lundi"
$ws.Range("B4").Value = "mercredi"
$ws.Range("B7").Value = "mercredi"

# Corresponding date serial numbers
$ws.Range("A2").Value = 46111
$ws.Range("A4").Value = 46113
$ws.Range("A7").Value = 46162
